$d = $word.ActiveDocument

# The Work Item table ("Table2" style) is the second table in the document.
$t = $d.Tables.Item(2)

# --- Table-level width / indent -------------------------------------------------
# tblW: 10946 dxa -> 10905 dxa  (547.3pt -> 545.25pt)
# tblInd: -740 dxa -> -695 dxa  (-37pt -> -34.75pt)
$t.PreferredWidth = 545.25
$t.Rows.LeftIndent = -34.75

# --- Column widths ---------------------------------------------------------------
# gridCol 1: 690 dxa -> 645 dxa (34.5pt -> 32.25pt)
# gridCol 8: 1136 dxa -> 1140 dxa (56.8pt -> 57pt)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $t.Cell($r, 1).Width = 32.25
    $t.Cell($r, 8).Width = 57.0
}

# --- Text edits --------------------------------------------------------------------
# Work item 4.4 ("Implement front end to support view organisational requirements"):
# "Assigned to" changes from Matthew to Jack.
$t.Cell(12, 5).Range.Text = "Jack"

# Work item 5.1 ("Completed Implementation of update project use case"):
# Outcome text tweak: "last iteration" -> "the last iteration".
$t.Cell(14, 3).Range.Text = "Carried over from the last iteration. Allows the update project use case to be fully realised"
